$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.490917
$ws.Range("H2").Value = 4.472751
$ws.Range("I2").Value = 0.02171654063137504
$ws.Range("J2").Value = 0.02171654063137504
$ws.Range("M2").Value = 10.34761366666667
$ws.Range("N2").Value = 31.042841
$ws.Range("O2").Value = 0.2299953477621856
$ws.Range("P2").Value = 0.2299953477621856
$ws.Range("Q2").Value = 15.42743312506566
$ws.Range("R2").Value = 138.846898125591
$ws.Range("S2").Value = 0.004994703314704734
$ws.Range("T2").Value = 0.004994703314704735
$ws.Range("G3").Value = 1.490917
$ws.Range("H3").Value = 4.472751
$ws.Range("I3").Value = 0.02171654063137504
$ws.Range("J3").Value = 0.02171654063137504
$ws.Range("O3").Value = 0.6794731949692173
$ws.Range("P3").Value = 0.6794731949692174
$ws.Range("Q3").Value = 45.57712744042633
$ws.Range("R3").Value = 410.1941469638369
$ws.Range("S3").Value = 0.01475580724647922
$ws.Range("T3").Value = 0.01475580724647922
$ws.Range("G4").Value = 1.490917
$ws.Range("H4").Value = 4.472751
$ws.Range("I4").Value = 0.02171654063137504
$ws.Range("J4").Value = 0.02171654063137504
$ws.Range("M4").Value = 4.073058666666666
$ws.Range("N4").Value = 12.219176
$ws.Range("O4").Value = 0.09053145726859702
$ws.Range("P4").Value = 0.09053145726859703
$ws.Range("Q4").Value = 6.072592408130665
$ws.Range("R4").Value = 54.65333167317599
$ws.Range("S4").Value = 0.00196603007019108
$ws.Range("T4").Value = 0.00196603007019108
$ws.Range("I5").Value = 0.5411744207383854
$ws.Range("J5").Value = 0.5411744207383854
$ws.Range("M5").Value = 10.34761366666667
$ws.Range("N5").Value = 31.042841
$ws.Range("O5").Value = 0.2299953477621856
$ws.Range("P5").Value = 0.2299953477621856
$ws.Range("Q5").Value = 384.4503747929099
$ws.Range("R5").Value = 3460.05337313619
$ws.Range("S5").Value = 0.1244675990977243
$ws.Range("T5").Value = 0.1244675990977243
$ws.Range("I6").Value = 0.5411744207383854
$ws.Range("J6").Value = 0.5411744207383854
$ws.Range("O6").Value = 0.6794731949692173
$ws.Range("P6").Value = 0.6794731949692174
$ws.Range("S6").Value = 0.3677135126947262
$ws.Range("T6").Value = 0.3677135126947262
$ws.Range("I7").Value = 0.5411744207383854
$ws.Range("J7").Value = 0.5411744207383854
$ws.Range("M7").Value = 4.073058666666666
$ws.Range("N7").Value = 12.219176
$ws.Range("O7").Value = 0.09053145726859702
$ws.Range("P7").Value = 0.09053145726859703
$ws.Range("Q7").Value = 151.32850736376
$ws.Range("R7").Value = 1361.95656627384
$ws.Range("S7").Value = 0.04899330894593488
$ws.Range("T7").Value = 0.04899330894593489
$ws.Range("G8").Value = 30.009075
$ws.Range("H8").Value = 90.027225
$ws.Range("I8").Value = 0.4371090386302395
$ws.Range("J8").Value = 0.4371090386302395
$ws.Range("M8").Value = 10.34761366666667
$ws.Range("N8").Value = 31.042841
$ws.Range("O8").Value = 0.2299953477621856
$ws.Range("P8").Value = 0.2299953477621856
$ws.Range("Q8").Value = 310.522314594025
$ws.Range("R8").Value = 2794.700831346225
$ws.Range("S8").Value = 0.1005330453497565
$ws.Range("T8").Value = 0.1005330453497566
$ws.Range("G9").Value = 30.009075
$ws.Range("H9").Value = 90.027225
$ws.Range("I9").Value = 0.4371090386302395
$ws.Range("J9").Value = 0.4371090386302395
$ws.Range("O9").Value = 0.6794731949692173
$ws.Range("P9").Value = 0.6794731949692174
$ws.Range("Q9").Value = 917.373291500675
$ws.Range("R9").Value = 8256.359623506076
$ws.Range("S9").Value = 0.2970038750280119
$ws.Range("T9").Value = 0.297003875028012
$ws.Range("G10").Value = 30.009075
$ws.Range("H10").Value = 90.027225
$ws.Range("I10").Value = 0.4371090386302395
$ws.Range("J10").Value = 0.4371090386302395
$ws.Range("M10").Value = 4.073058666666666
$ws.Range("N10").Value = 12.219176
$ws.Range("O10").Value = 0.09053145726859702
$ws.Range("P10").Value = 0.09053145726859703
$ws.Range("Q10").Value = 122.2287230074
$ws.Range("R10").Value = 1100.0585070666
$ws.Range("S10").Value = 0.03957211825247105
$ws.Range("T10").Value = 0.03957211825247105
